$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append after the existing A1:B25 table.
$ws.Range("A26").Value = 45993
$ws.Range("B26").Value = 2

$ws.Range("A27").Value = 45994
$ws.Range("B27").Value = 1

# Match the date formatting used by the rest of column A (e.g. A25) without
# introducing a new style entry - copy formats only, like Excel's "paste
# special > formats" after filling the row.
$ws.Range("A25").Copy() | Out-Null
$ws.Range("A26:A27").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Reflect the updated scroll position / selection, as in the saved file.
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("A26:B27").Select()
